$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing row 54 (and below) down to row 55
$ws.Rows.Item(54).Insert()

# Fill in the new row 54 with data
$ws.Cells.Item(54, 1).Value = 8
$ws.Cells.Item(54, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value = 44568
$ws.Cells.Item(54, 5).Value = 4
$ws.Cells.Item(54, 6).Value = 100112030
$ws.Cells.Item(54, 7).Value = "Poroto granado"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 700
$ws.Cells.Item(54, 11).Value = 29000
$ws.Cells.Item(54, 12).Value = 30000
$ws.Cells.Item(54, 13).Value = 29500
$ws.Cells.Item(54, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(54, 16).Value = 1180
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = "Hortaliza"
